$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.421.59"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.14%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.908.95"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.59%  "

# Row 4
$ws.Range("E4").Value = "  +0.32%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "238.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.64%  "

# Row 6
$ws.Range("E6").Value = "  +0.20%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4766"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.77%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2855"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.67%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06666"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.29%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.79"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.99%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "102.40"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.30%  "

# Row 12
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.921.00"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.08%  "

# Row 13
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07713"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.10%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.214"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.85%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6748"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.85%  "

# Row 16
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "30.447.89"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.11%  "

# Row 17
$ws.Range("B17").Value = "BitcoinCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "259.01"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -8.49%  "

# Row 18
$ws.Range("E18").Value = "  +0.14%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007469"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.56%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.68"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.37%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.410"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.72%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.005"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.42%  "

# Row 23
$ws.Range("B23").Value = "BitDAO"
$ws.Range("C23").Value = "https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.4583"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -8.40%  "

# Row 24
$ws.Range("B24").Value = "Chainlink"
$ws.Range("C24").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.288"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.40%  "

# Row 25
$ws.Range("B25").Value = "Cosmos"
$ws.Range("C25").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.440"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.11%  "

# Row 26
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "162.66"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.49%  "

# Row 27
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.87"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.84%  "

# Row 28
$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.071"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.98%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.382"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.25%  "

# Row 30
$ws.Range("B30").Value = "Stellar"
$ws.Range("C30").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.1007"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.56%  "

# Row 31
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.589"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.66%  "

# Row 32
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.512"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.25%  "

# Row 33
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.205"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.60%  "

# Row 34
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04784"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.25%  "

# Row 35
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7274"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.17%  "

# Row 36
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.113"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.11%  "

# Row 37
$ws.Range("B37").Value = "Frax"
$ws.Range("C37").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.002"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.16%  "

# Row 38
$ws.Range("B38").Value = "HuobiToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.712"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.82%  "

# Row 39
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01921"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.85%  "

# Row 40
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.632"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.65%  "

# Row 41
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.265"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.26%  "

# Row 42
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "74.20"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.82%  "

# Row 43
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.991"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.54%  "

# Row 44
$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "106.67"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.49%  "

# Row 45
$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8599"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.15%  "

# Row 46
$ws.Range("B46").Value = "TheSandbox"
$ws.Range("C46").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4250"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.17%  "

# Row 47
$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.002"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.09%  "

# Row 48
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.008.17"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.59%  "

# Row 49
$ws.Range("B49").Value = "Aptos"
$ws.Range("C49").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.439"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -9.17%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "35.04"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.60%  "

# Row 51
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.1193"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.41%  "
